# Fix a typo in filter slide
#
# The presentation's "filter" slide (the one whose shapes contain the
# ".filter(item => item % 2 === 0)" TypeScript snippet and the matching
# "out.add(item);" pseudo-code line) has two small text corrections:
#
#   1. `const target: string[] = ...`  ->  `const target: number[] = ...`
#      (the filtered array actually holds numbers, not strings)
#   2. `       out.add(item);`          ->  `       out.add(item); // clone!!!`
#      (a reminder comment that items are cloned into the output collection)
#
# Both edits are applied by locating the exact paragraph that holds the
# text and rewriting just that paragraph. To keep each paragraph's run
# formatting (font, size, color, ...) intact we first overwrite the
# paragraph with a placeholder string that shares no characters with
# either the old or new text, and only then assign the real text. This
# avoids the host's text-diff logic splitting the paragraph into extra
# runs because of a shared prefix/suffix between old and new text, which
# would otherwise needlessly fragment a single run into several
# identically formatted ones.
#
# NOTE: this host's function calls only reliably bind positional
# parameters (named parameters such as "-Slide $x" and parameter default
# values are not honored), so every helper below takes plain positional
# arguments and every call site passes all of them explicitly.

function Set-ParagraphText {
    param($TextRange, $OldText, $NewText, $MaxParagraphs)

    # A paragraph's .Text includes a trailing CR (chr 13) for every
    # paragraph except the very last one in the text frame when *read*,
    # but writing that CR back is taken literally (it ends up as a
    # run containing a newline character) instead of being consumed as
    # the paragraph boundary marker. So: compare with it stripped, but
    # never write it back - omitting it on write keeps the paragraph
    # boundary intact anyway.
    $cr = [char]13

    for ($i = 1; $i -le $MaxParagraphs; $i++) {
        $para = $TextRange.Paragraphs($i, 1)
        $plain = $para.Text.TrimEnd($cr)

        if ($plain -eq $OldText) {
            # Neutral placeholder with no characters in common with the
            # old/new text so the host can't find a shared prefix/suffix
            # and therefore rewrites the paragraph as a single clean run.
            $para.Text = "@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@"
            $para2 = $TextRange.Paragraphs($i, 1)
            $para2.Text = $NewText
            return $true
        }
    }
    return $false
}

function Find-ShapeWithText {
    param($Slide, $Needle)

    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shape = $Slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.Contains($Needle)) {
                return $shape
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

for ($slideIndex = 1; $slideIndex -le $p.Slides.Count; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)

    # --- Edit 1: filter() snippet's array type annotation ---------------
    $shape1 = Find-ShapeWithText $slide "const target: string[] = [1, 2, 3, 4, 5].filter(item"
    if ($null -ne $shape1) {
        Set-ParagraphText $shape1.TextFrame.TextRange "const target: string[] = [1, 2, 3, 4, 5].filter(item => item % 2 === 0); // [2, 4]" "const target: number[] = [1, 2, 3, 4, 5].filter(item => item % 2 === 0); // [2, 4]" 200
    }

    # --- Edit 2: pseudo-code "out.add(item);" line -----------------------
    $shape2 = Find-ShapeWithText $slide "out.add(item);"
    if ($null -ne $shape2) {
        Set-ParagraphText $shape2.TextFrame.TextRange "       out.add(item);" "       out.add(item); // clone!!!" 200
    }
}
